# Update generated "想去人数" (want-to-go count) and "最低票价" (lowest price)
# figures on both the "展览" and "全部类型" worksheets (they mirror the
# same data), matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 7336
    $ws.Range("G3").Value = 60

    $ws.Range("F4").Value = 5686

    $ws.Range("F6").Value = 178

    $ws.Range("F13").Value = 61

    $ws.Range("F15").Value = 381

    $ws.Range("F17").Value = 16

    $ws.Range("F20").Value = 55
}
